$wb = $excel.ActiveWorkbook

# --- Metadata sheet text updates (bump term version) ---
$ws = $wb.Worksheets.Item("Metadata")

# Locate the "Version" and "Date" property rows, then update the values
# next to them in column B.
$versionLabel = $ws.Range("A1:A20").Find("Version")
$ws.Cells.Item($versionLabel.Row, $versionLabel.Column + 1).Value = "1.1.0"

$dateLabel = $ws.Range("A1:A20").Find("Date")
$ws.Cells.Item($dateLabel.Row, $dateLabel.Column + 1).Value = "2023-07-10T23:08:03+02:00"

# --- Style updates: make sure the wrap-text alignment that is defined on
#     the header-row style and the regular-row style is actually flagged
#     as applied (applyAlignment), so the existing wrapText/vertical-top
#     alignment settings take effect on every sheet. ---
foreach ($sheet in $wb.Worksheets) {
    $sheet.UsedRange.WrapText = $true
}
